# Auto-generated COM-interop script to add AL column (sheet "data") and AK column (sheet "pocetR")
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "data" ---
$ws1 = $wb.Worksheets.Item("data")

# Header cell AL1: copy formatting from AK1 (same style as other date headers), then set value/date text
$ws1.Range("AK1").Copy() | Out-Null
$ws1.Range("AL1").PasteSpecial(-4122) | Out-Null
$ws1.Range("AL1").Value = "28. 12. 2021"

# Data rows 2-61: new AL value for the "28. 12. 2021" wave
$ws1.Range("AL2").Value = 0.53
$ws1.Range("AL3").Value = 0.29
$ws1.Range("AL4").Value = 0.18
$ws1.Range("AL5").Value = 0.36
$ws1.Range("AL6").Value = 0.18
$ws1.Range("AL7").Value = 0.46
$ws1.Range("AL8").Value = 0.54
$ws1.Range("AL9").Value = 0.3
$ws1.Range("AL10").Value = 0.16
$ws1.Range("AL11").Value = 0.54
$ws1.Range("AL12").Value = 0.31
$ws1.Range("AL13").Value = 0.15
$ws1.Range("AL14").Value = 0.5
$ws1.Range("AL15").Value = 0.24
$ws1.Range("AL16").Value = 0.26
$ws1.Range("AL17").Value = 0.55
$ws1.Range("AL18").Value = 0.3
$ws1.Range("AL19").Value = 0.15
$ws1.Range("AL20").Value = 0.49
$ws1.Range("AL21").Value = 0.27
$ws1.Range("AL22").Value = 0.24
$ws1.Range("AL23").Value = 0.41
$ws1.Range("AL24").Value = 0.25
$ws1.Range("AL25").Value = 0.34
$ws1.Range("AL26").Value = 0.43
$ws1.Range("AL27").Value = 0.34
$ws1.Range("AL28").Value = 0.23
$ws1.Range("AL29").Value = 0.57
$ws1.Range("AL30").Value = 0.28
$ws1.Range("AL31").Value = 0.15
$ws1.Range("AL32").Value = 0.61
$ws1.Range("AL33").Value = 0.24
$ws1.Range("AL34").Value = 0.15
$ws1.Range("AL35").Value = 0.44
$ws1.Range("AL36").Value = 0.3
$ws1.Range("AL37").Value = 0.26
$ws1.Range("AL38").Value = 0.39
$ws1.Range("AL39").Value = 0.39
$ws1.Range("AL40").Value = 0.22
$ws1.Range("AL41").Value = 0.59
$ws1.Range("AL42").Value = 0.24
$ws1.Range("AL43").Value = 0.17
$ws1.Range("AL44").Value = 0.67
$ws1.Range("AL45").Value = 0.22
$ws1.Range("AL46").Value = 0.11
$ws1.Range("AL47").Value = 0.51
$ws1.Range("AL48").Value = 0.3
$ws1.Range("AL49").Value = 0.19
$ws1.Range("AL50").Value = 0.68
$ws1.Range("AL51").Value = 0.23
$ws1.Range("AL52").Value = 0.09
$ws1.Range("AL53").Value = 0.44
$ws1.Range("AL54").Value = 0.33
$ws1.Range("AL55").Value = 0.23
$ws1.Range("AL56").Value = 0.59
$ws1.Range("AL57").Value = 0.23
$ws1.Range("AL58").Value = 0.18
$ws1.Range("AL59").Value = 0.67
$ws1.Range("AL60").Value = 0.26
$ws1.Range("AL61").Value = 0.07

# Row 62 caption: bump the "aktualizace" (update) date in the footer text
$ws1.Range("A62").Value = "Život během pandemie, Obavy ze ztráty práce, % respondentů celkově a ve skupinách, aktualizace 6. 1. 2022"

# --- Sheet 2: "pocetR" ---
$ws2 = $wb.Worksheets.Item("pocetR")

# Header cell AK1: copy formatting from AJ1, then set value
$ws2.Range("AJ1").Copy() | Out-Null
$ws2.Range("AK1").PasteSpecial(-4122) | Out-Null
$ws2.Range("AK1").Value = "28. 12. 2021"

# Data rows 2-24: new AK sample-size values for the "28. 12. 2021" wave
$ws2.Range("AK2").Value = 994
$ws2.Range("AK3").Value = 74
$ws2.Range("AK4").Value = 920
$ws2.Range("AK5").Value = 777
$ws2.Range("AK6").Value = 141
$ws2.Range("AK7").Value = 6
$ws2.Range("AK8").Value = 70
$ws2.Range("AK9").Value = 744
$ws2.Range("AK10").Value = 128
$ws2.Range("AK11").Value = 65
$ws2.Range("AK12").Value = 57
$ws2.Range("AK13").Value = 381
$ws2.Range("AK14").Value = 380
$ws2.Range("AK15").Value = 233
$ws2.Range("AK16").Value = 124
$ws2.Range("AK17").Value = 313
$ws2.Range("AK18").Value = 316
$ws2.Range("AK19").Value = 151
$ws2.Range("AK20").Value = 273
$ws2.Range("AK21").Value = 87
$ws2.Range("AK22").Value = 275
$ws2.Range("AK23").Value = 149
$ws2.Range("AK24").Value = 79

# Row 25 caption + trailing blank cell (extends the used range like the empty
# inlineStr cells already present at B25:AJ25) to match used range extension
$ws2.Range("A25").Value = "Život během pandemie, Obavy ze ztráty práce, velikost dotázaného souboru celkově a ve skupinách, aktualizace 6. 1. 2022"
$ws2.Range("AJ25").Copy() | Out-Null
$ws2.Range("AK25").PasteSpecial(-4122) | Out-Null

